# Auto-generated script to update Leve profit calculation columns (H-N)
# across multiple worksheets, per scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 133.07692
$ws.Range("I33").Value = 86.19047500000001
$ws.Range("J33").Value = 330
$ws.Range("K33").Value = 86.19047500000001
$ws.Range("L33").Value = 330
$ws.Range("M33").Value = 142.809525
$ws.Range("N33").Value = -788

$ws.Range("H37").Value = 500
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 500
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 1500
$ws.Range("N37").Value = -1752

$ws.Range("H39").Value = 51.214287
$ws.Range("I39").Value = 31.222221
$ws.Range("J39").Value = 87.2
$ws.Range("K39").Value = 93.666663
$ws.Range("L39").Value = 261.6
$ws.Range("M39").Value = 202.333337

$ws.Range("H97").Value = 3158.3333
$ws.Range("I97").Value = 8420
$ws.Range("J97").Value = 2106
$ws.Range("K97").Value = 25260
$ws.Range("L97").Value = 6318
$ws.Range("M97").Value = -24764
$ws.Range("N97").Value = -7310

$ws.Range("H106").Value = 35730940
$ws.Range("I106").Value = 38476360
$ws.Range("J106").Value = 40503
$ws.Range("K106").Value = 38476360
$ws.Range("L106").Value = 40503
$ws.Range("M106").Value = -38475729

$ws.Range("H121").Value = 6949.6
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 6949.6
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 20848.8
$ws.Range("N121").Value = -24342.8

$ws.Range("H137").Value = 2160
$ws.Range("I137").Value = 2046
$ws.Range("J137").Value = 2274
$ws.Range("K137").Value = 6138
$ws.Range("L137").Value = 6822
$ws.Range("M137").Value = -3588
$ws.Range("N137").Value = -11922

$ws.Range("H138").Value = 2097.6
$ws.Range("I138").Value = 1951.6666
$ws.Range("J138").Value = 2150.6667
$ws.Range("K138").Value = 5854.9998
$ws.Range("L138").Value = 6452.000100000001
$ws.Range("M138").Value = -714.9997999999996
$ws.Range("N138").Value = -16732.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4719.3335
$ws.Range("I32").Value = 3311.6428
$ws.Range("J32").Value = 12602.4
$ws.Range("K32").Value = 3311.6428
$ws.Range("L32").Value = 12602.4
$ws.Range("M32").Value = -3024.6428

$ws.Range("H61").Value = 6165.8335
$ws.Range("I61").Value = 6165.8335
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6165.8335
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5953.8335

$ws.Range("H76").Value = 26632.666
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 26632.666
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 26632.666
$ws.Range("N76").Value = -27308.666

$ws.Range("H79").Value = 26632.666
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 26632.666
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 26632.666
$ws.Range("N79").Value = -28972.666

$ws.Range("H136").Value = 6165.8335
$ws.Range("I136").Value = 6165.8335
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 18497.5005
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -15947.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 7790
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 7790
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 7790
$ws.Range("N103").Value = -10134

$ws.Range("H134").Value = 1347
$ws.Range("I134").Value = 1396.7778
$ws.Range("J134").Value = 899
$ws.Range("K134").Value = 4190.3334
$ws.Range("L134").Value = 2697
$ws.Range("M134").Value = -1655.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3429.8965
$ws.Range("I31").Value = 2300.1
$ws.Range("J31").Value = 5940.5557
$ws.Range("K31").Value = 2300.1
$ws.Range("L31").Value = 5940.5557
$ws.Range("M31").Value = -2005.1

$ws.Range("H34").Value = 3429.8965
$ws.Range("I34").Value = 2300.1
$ws.Range("J34").Value = 5940.5557
$ws.Range("K34").Value = 2300.1
$ws.Range("L34").Value = 5940.5557
$ws.Range("M34").Value = -2098.1

$ws.Range("H60").Value = 12563.1
$ws.Range("I60").Value = 10626.333
$ws.Range("J60").Value = 29994
$ws.Range("K60").Value = 10626.333
$ws.Range("L60").Value = 29994
$ws.Range("M60").Value = -10115.333

$ws.Range("H107").Value = 13514421
$ws.Range("I107").Value = 22727780
$ws.Range("J107").Value = 1494.2
$ws.Range("K107").Value = 22727780
$ws.Range("L107").Value = 1494.2
$ws.Range("M107").Value = -22725860
$ws.Range("N107").Value = -5334.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 262.07693
$ws.Range("I2").Value = 37.22222
$ws.Range("J2").Value = 768
$ws.Range("K2").Value = 223.33332
$ws.Range("L2").Value = 4608
$ws.Range("M2").Value = -110.33332
$ws.Range("N2").Value = -4834

$ws.Range("H8").Value = 193
$ws.Range("I8").Value = 193
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 579
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -440

$ws.Range("H12").Value = 142.72223
$ws.Range("I12").Value = 154.44444
$ws.Range("J12").Value = 131
$ws.Range("K12").Value = 463.33332
$ws.Range("L12").Value = 393
$ws.Range("M12").Value = -290.33332
$ws.Range("N12").Value = -739

$ws.Range("H14").Value = 373.54544
$ws.Range("I14").Value = 373.54544
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1120.63632
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -947.6363200000001

$ws.Range("H33").Value = 84235.75
$ws.Range("I33").Value = 1189.7778
$ws.Range("J33").Value = 333373.66
$ws.Range("K33").Value = 7138.666800000001
$ws.Range("L33").Value = 2000241.96
$ws.Range("M33").Value = -6855.666800000001
$ws.Range("N33").Value = -2000807.96

$ws.Range("H68").Value = 644.61536
$ws.Range("I68").Value = 603.44446
$ws.Range("J68").Value = 737.25
$ws.Range("K68").Value = 1810.33338
$ws.Range("L68").Value = 2211.75
$ws.Range("M68").Value = -999.33338
$ws.Range("N68").Value = -3833.75

$ws.Range("H71").Value = 644.61536
$ws.Range("I71").Value = 603.44446
$ws.Range("J71").Value = 737.25
$ws.Range("K71").Value = 5431.00014
$ws.Range("L71").Value = 6635.25
$ws.Range("M71").Value = -1375.00014
$ws.Range("N71").Value = -14747.25

$ws.Range("H122").Value = 562.4286
$ws.Range("I122").Value = 291.5
$ws.Range("J122").Value = 670.8
$ws.Range("K122").Value = 2623.5
$ws.Range("L122").Value = 6037.2
$ws.Range("M122").Value = -173.5
$ws.Range("N122").Value = -10937.2

$ws.Range("H140").Value = 5424
$ws.Range("I140").Value = 2844.625
$ws.Range("J140").Value = 6799.6665
$ws.Range("K140").Value = 8533.875
$ws.Range("L140").Value = 20398.9995
$ws.Range("M140").Value = -3353.875
$ws.Range("N140").Value = -30758.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 173388.67
$ws.Range("I7").Value = 505000
$ws.Range("J7").Value = 7583
$ws.Range("K7").Value = 505000
$ws.Range("L7").Value = 7583
$ws.Range("M7").Value = -504888
$ws.Range("N7").Value = -7807

$ws.Range("H8").Value = 173388.67
$ws.Range("I8").Value = 505000
$ws.Range("J8").Value = 7583
$ws.Range("K8").Value = 505000
$ws.Range("L8").Value = 7583
$ws.Range("M8").Value = -504861
$ws.Range("N8").Value = -7861

$ws.Range("I122").Value = 2149.5
$ws.Range("J122").Value = 128047.625
$ws.Range("K122").Value = 6448.5
$ws.Range("L122").Value = 384142.875
$ws.Range("M122").Value = -3998.5
$ws.Range("N122").Value = -389042.875

$ws.Range("H126").Value = 1999.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 4800
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -3774

$ws.Range("H21").Value = 6500
$ws.Range("I21").Value = 2000
$ws.Range("J21").Value = 7142.857
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 7142.857
$ws.Range("M21").Value = -1826

$ws.Range("H22").Value = 1167
$ws.Range("I22").Value = 1250.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1250.5
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -955.5
$ws.Range("N22").Value = -1590

$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 10000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10686

$ws.Range("H27").Value = 1167
$ws.Range("I27").Value = 1250.5
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1250.5
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -1143.5
$ws.Range("N27").Value = -1214

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = ""

$ws.Range("H68").Value = 4589.4
$ws.Range("I68").Value = 4314.6665
$ws.Range("J68").Value = 5001.5
$ws.Range("K68").Value = 4314.6665
$ws.Range("L68").Value = 5001.5
$ws.Range("M68").Value = -3565.6665

$ws.Range("H71").Value = 4589.4
$ws.Range("I71").Value = 4314.6665
$ws.Range("J71").Value = 5001.5
$ws.Range("K71").Value = 21573.3325
$ws.Range("L71").Value = 25007.5
$ws.Range("M71").Value = -17829.3325

$ws.Range("H82").Value = 57872.445
$ws.Range("I82").Value = 2600.4285
$ws.Range("J82").Value = 251324.5
$ws.Range("K82").Value = 2600.4285
$ws.Range("L82").Value = 251324.5
$ws.Range("M82").Value = -2239.4285
$ws.Range("N82").Value = -252046.5

$ws.Range("H85").Value = 57872.445
$ws.Range("I85").Value = 2600.4285
$ws.Range("J85").Value = 251324.5
$ws.Range("K85").Value = 2600.4285
$ws.Range("L85").Value = 251324.5
$ws.Range("M85").Value = -1352.4285
$ws.Range("N85").Value = -253820.5

$ws.Range("H100").Value = 6399.8
$ws.Range("I100").Value = 6399.8
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 6399.8
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -5858.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8074
$ws.Range("I81").Value = 5105.7144
$ws.Range("J81").Value = 15000
$ws.Range("K81").Value = 10211.4288
$ws.Range("L81").Value = 30000
$ws.Range("M81").Value = -9150.4288

$ws.Range("H84").Value = 8074
$ws.Range("I84").Value = 5105.7144
$ws.Range("J84").Value = 15000
$ws.Range("K84").Value = 51057.144
$ws.Range("L84").Value = 150000
$ws.Range("M84").Value = -45753.144

$ws.Range("H113").Value = 1949.25
$ws.Range("I113").Value = 765.6667
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 2297.0001
$ws.Range("L113").Value = 16500
$ws.Range("M113").Value = -127.0001000000002
$ws.Range("N113").Value = -20840
